# IGCC_Netting_Flows_Historical.xlsx update ("Updating the model for Horeco")
#
# The model rolls its 96-quarter-per-day timestamp table forward by one
# calendar day: every "Timestamp (CET)" value in column A is bumped +1 day,
# and the literal "Lookup" text in column E (which encodes the same date as
# DD.MM.YYYY immediately followed by the Quarter number from column D) is
# regenerated to stay in sync. Columns B, C and D are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Range("A$r")
    $dCell = $ws.Range("D$r")
    $eCell = $ws.Range("E$r")

    $oldSerial = $aCell.Value2
    $newSerial = $oldSerial + 1
    $aCell.Value = $newSerial

    $quarter = $dCell.Value2
    $newDate = [datetime]::FromOADate($newSerial)
    $eCell.Value = $newDate.ToString("dd.MM.yyyy") + [string]$quarter
}
